$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.016.10"
$ws.Range("E2").Value = "  +2.56%  "

$ws.Range("D3").Value = "3.737.65"
$ws.Range("E3").Value = "  +1.24%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'601.66"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").Value = "'167.92"
$ws.Range("E6").Value = "  +1.18%  "

$ws.Range("D7").Value = "3.734.99"
$ws.Range("E7").Value = "  +1.12%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +0.30%  "

$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("D11").Value = "'6.43"
$ws.Range("E11").Value = "  +3.78%  "

$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").Value = "'37.96"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("E14").Value = "  +1.70%  "

$ws.Range("D15").Value = "4.361.64"
$ws.Range("E15").Value = "  +1.31%  "

$ws.Range("D16").Value = "3.738.68"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").Value = "69.050.30"
$ws.Range("E17").Value = "  +2.67%  "

$ws.Range("D18").Value = "'7.29"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("E19").Value = "  -1.17%  "

$ws.Range("D20").Value = "'17.08"
$ws.Range("E20").Value = "  -0.88%  "

$ws.Range("D21").Value = "'10.79"
$ws.Range("E21").Value = "  +17.39%  "

$ws.Range("D22").Value = "'492.46"
$ws.Range("E22").Value = "  +1.46%  "

$ws.Range("D23").Value = "'0.724"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("E24").Value = "  +6.92%  "

$ws.Range("D25").Value = "'84.71"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("E26").Value = "  +0.76%  "

$ws.Range("D27").Value = "'12.29"
$ws.Range("E27").Value = "  +0.65%  "

$ws.Range("D28").Value = "'10.12"
$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  +2.40%  "

$ws.Range("D31").Value = "'2.50"
$ws.Range("E31").Value = "  +6.56%  "

$ws.Range("E32").Value = "  +4.79%  "

$ws.Range("D33").Value = "'31.50"
$ws.Range("E33").Value = "  +0.92%  "

$ws.Range("D34").Value = "3.882.37"
$ws.Range("E34").Value = "  +1.44%  "

$ws.Range("D35").Value = "'0.109"
$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("D36").Value = "3.671.50"
$ws.Range("E36").Value = "  +1.08%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  +1.77%  "

$ws.Range("D39").Value = "'5.85"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("D41").Value = "'0.324"
$ws.Range("E41").Value = "  +0.95%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.96"
$ws.Range("E42").Value = "  +5.69%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'432.65"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("D44").Value = "'48.52"
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("D45").Value = "'1.98"
$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("D46").Value = "'8.46"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").Value = "'40.09"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("D49").Value = "'141.78"
$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("D50").Value = "2.773.91"
$ws.Range("E50").Value = "  +1.16%  "

$ws.Range("D51").Value = "'0.0353"
$ws.Range("E51").Value = "  +1.04%  "
